# "saving latest uart stuff"
# Update the desired Tloop (ms) input on Sheet1 (cell C9), which drives all
# of the downstream formulas/results on both Sheet1 and Sheet2, and move the
# active cell selection to C9 (matching where the user left off editing).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update the user-entered "desired Tloop (ms)" value; everything else on
# both sheets recalculates automatically from this single input.
$ws1.Range("C9").Value = 500

# Reflect the new active cell/selection on Sheet1.
$ws1.Activate()
$ws1.Range("C9").Select()
